$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-05-17"

# Update the header label for the current-year column (I1) to match new date
$ws.Range("I1").Value = "2022 (through 05-17)"

# Update the May total for the 2022 column (I6) with the new count
$ws.Range("I6").Value = 62

# Update the grand total for the 2022 column (I14) with the new count
$ws.Range("I14").Value = 614
